$d = $word.ActiveDocument

# --- 1. Add the "_GoBack" bookmark (zero-length) right at the very start of
#        the document body, matching what Word leaves behind after an edit
#        session. Adding a bookmark at absolute position 0 directly can make
#        its end expand across the first paragraph, so we use a temporary
#        placeholder character, anchor the bookmark after it, then remove the
#        placeholder; the bookmark collapses back to position 0 cleanly.
$pad = $d.Range(0, 0)
$pad.InsertBefore("Z")
$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range(0, 1).Delete()

# --- 2. Remove the big block of duplicated "Q1..Q40" interview-question
#        paragraphs (with tab separators) that had been appended after the
#        original "Q11. What is a VLAN ..." paragraph, restoring the document
#        to its pre-append state. This also removes the trailing empty
#        paragraph that was left at the very end of the body.
$target = "What is the purpose of a computer network, and how does it facilitate communication between devices?"
$startPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match [regex]::Escape("Q1") -and $p.Range.Text -match [regex]::Escape($target)) {
        $startPara = $i
        break
    }
}

if ($startPara -ne $null) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $first = $d.Paragraphs.Item($startPara)
    $r = $d.Range($first.Range.Start, $lastPara.Range.End)
    $r.Delete()
}
